# Adds the new "Text jaText has been added " text box to slide 1 and
# turns on the (empty) PowerPoint-2012 slide-guide-list presentation extension,
# matching the authored diff.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- New text box on slide 1 -------------------------------------------------
# Position/size below are the EMU values from the diff, converted to points
# (1 pt = 12700 EMU) since Shapes.AddTextbox takes points.
$left   = 12435840 / 12700
$top    = 4941332  / 12700
$width  = 10698480 / 12700
$height = 369332   / 12700

$shp = $s.Shapes.AddTextbox(1, $left, $top, $width, $height)

$shp.Rotation = 180
$shp.VerticalFlip = $true

$shp.TextFrame.WordWrap = $true
$shp.TextFrame.AutoSize = 1

$shp.Fill.Visible = $false

$tr = $shp.TextFrame.TextRange
$tr.Text = "Text "
$tr2 = $tr.InsertAfter("jaText")
$tr3 = $tr2.InsertAfter(" has ")
$tr4 = $tr3.InsertAfter("been added ")
